# Apply weekly refresh: a new observation date (2025-01-27, serial 45684)
# became available, which:
#   1) fills in the previously-unknown "actual" values that were predicted
#      in earlier weeks (diagonal cascade in the "Valeurs reelles" sheet), and
#   2) appends a brand-new row for the new date on both sheets, with the
#      freshly observed value (sheet 1) / freshly computed predictions (sheet 2).

$wb = $excel.ActiveWorkbook

$wsValeurs = $wb.Worksheets.Item(1)   # "Valeurs réelles"
$wsPred    = $wb.Worksheets.Item(2)   # "Prédictions"

# Append the new row 26 (date 2025-01-27). Copy row 25 down first - while it
# still holds its original (un-cascaded) values - so the untouched trailing
# cells keep the same "blank" shape/formatting as the rest of the sheet, then
# overwrite the cells that actually have new data.
$wsValeurs.Range("A25:E25").Copy($wsValeurs.Range("A26:E26"))
$wsValeurs.Range("A26").Value = 45684
$wsValeurs.Range("B26").Value = 1.210000038146973

# The new actual value (1.210000038146973) resolves cells that were previously
# unknown (shown blank) for dates whose S+1 / S+2 / S+3 horizon is this new date.
$wsValeurs.Range("E23").Value = 1.210000038146973
$wsValeurs.Range("D24").Value = 1.210000038146973
$wsValeurs.Range("C25").Value = 1.210000038146973

# --- Sheet 2: "Prédictions" ----------------------------------------------------
# Append the new row 26 with the freshly computed predictions.
$wsPred.Range("A25:D25").Copy($wsPred.Range("A26:D26"))
$wsPred.Range("A26").Value = 45684
$wsPred.Range("B26").Value = 0.969795823097229
$wsPred.Range("C26").Value = 0.9494584798812866
$wsPred.Range("D26").Value = 0.9467955827713013
